$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.435.20'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '1.636.19'
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '303.81'
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").Value = '0.3787'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '51.65'
$ws.Range("E8").Value = '  -0.96%  '
$ws.Range("D9").Value = '0.3622'
$ws.Range("E9").Value = '  -1.34%  '
$ws.Range("D10").Value = '0.08183'
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").Value = '1.229'
$ws.Range("E11").Value = '  -3.41%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '22.40'
$ws.Range("E13").Value = '  -3.36%  '
$ws.Range("D14").Value = '6.472'
$ws.Range("E14").Value = '  -3.61%  '
$ws.Range("D15").Value = '7.374'
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("E16").Value = '  -2.98%  '
$ws.Range("D17").Value = '1.629.96'
$ws.Range("E17").Value = '  -1.80%  '
$ws.Range("D18").Value = '95.12'
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("D19").Value = '0.06945'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").Value = '6.584'
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("D21").Value = '17.47'
$ws.Range("E21").Value = '  -5.48%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '12.56'
$ws.Range("E23").Value = '  -3.06%  '
$ws.Range("D24").Value = '23.407.70'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("E25").Value = '  +4.11%  '
$ws.Range("D26").Value = '3.053'
$ws.Range("E26").Value = '  -3.21%  '
$ws.Range("D27").Value = '21.14'
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("D28").Value = '150.92'
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("D29").Value = '5.274'
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("D30").Value = '133.35'
$ws.Range("E30").Value = '  -2.84%  '
$ws.Range("D31").Value = '1.810.87'
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = '2.166'
$ws.Range("E32").Value = '  -7.10%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.605'
$ws.Range("E33").Value = '  -4.64%  '
$ws.Range("D34").Value = '1.044'
$ws.Range("E34").Value = '  +6.59%  '
$ws.Range("D35").Value = '11.23'
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("D36").Value = '0.02749'
$ws.Range("E36").Value = '  -4.71%  '
$ws.Range("D37").Value = '0.08778'
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("D38").Value = '0.2493'
$ws.Range("E38").Value = '  -3.25%  '
$ws.Range("D39").Value = '0.07105'
$ws.Range("E39").Value = '  -4.09%  '
$ws.Range("D40").Value = '6.021'
$ws.Range("E40").Value = '  -5.99%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.7001'
$ws.Range("E41").Value = '  -3.36%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.340'
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("D43").Value = '15.89'
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("E44").Value = '  -4.06%  '
$ws.Range("D45").Value = '0.6501'
$ws.Range("E45").Value = '  -2.62%  '
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("E47").Value = '  -4.49%  '
$ws.Range("D48").Value = '3.969'
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("D49").Value = '0.07980'
$ws.Range("E49").Value = '  -0.74%  '
$ws.Range("D50").Value = '127.22'
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("D51").Value = '1.190'
$ws.Range("E51").Value = '  -3.33%  '
